# Updates the Turkey Super Lig 2023-2024 results sheet:
#   - Rows 49/50, 51/52, 55/56 and 94/95 had their match data (columns F:V)
#     swapped with each other.
#   - Rows 97/98/99 had their match data (columns F:V) cyclically rotated
#     (new 97 = old 99, new 98 = old 97, new 99 = old 98).
#   - Two new match rows (109 and 110) were appended at the end of the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that hold the per-match payload (team names, scores, odds, dates, url).
$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

function Get-RowValues($row) {
    $vals = @{}
    foreach ($c in $cols) {
        # NOTE: bare ".Value" (no parens) misbehaves as an rvalue in this
        # COM shim (returns a property-descriptor string instead of the
        # actual data) - ".Value2" is the reliable getter.
        $vals[$c] = $ws.Range("$c$row").Value2
    }
    return $vals
}

function Set-RowValues($row, $vals) {
    foreach ($c in $cols) {
        $ws.Range("$c$row").Value = $vals[$c]
    }
}

# --- Swap pairs: (49,50) (51,52) (55,56) (94,95) ---
$pairs = @(
    @(49, 50),
    @(51, 52),
    @(55, 56),
    @(94, 95)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $v1 = Get-RowValues $r1
    $v2 = Get-RowValues $r2
    Set-RowValues $r1 $v2
    Set-RowValues $r2 $v1
}

# --- Rotation: new97 = old99, new98 = old97, new99 = old98 ---
$v97 = Get-RowValues 97
$v98 = Get-RowValues 98
$v99 = Get-RowValues 99

Set-RowValues 97 $v99
Set-RowValues 98 $v97
Set-RowValues 99 $v98

# --- Append two brand-new rows (109, 110) ---

# Column A (bold, bordered, centered) and column E (date/time number format)
# carry explicit cell styles on every data row; replicate them from the last
# existing row (108) onto the two new rows before filling in values.
$ws.Range("A108").Copy() | Out-Null
$ws.Range("A109").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A110").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("E108").Copy() | Out-Null
$ws.Range("E109").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("E110").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("A109").Value = 108
$ws.Range("B109").Value = "turkey"
$ws.Range("C109").Value = "super-lig"
$ws.Range("D109").Value = "2023-2024"
$ws.Range("E109").Value = 45236.75
$ws.Range("F109").Value = "Gaziantep"
$ws.Range("G109").Value = 2
$ws.Range("H109").Value = "Rizespor"
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 2.32
$ws.Range("K109").Value = "30/10/2023 18:12"
$ws.Range("L109").Value = 2.2
$ws.Range("M109").Value = "06/11/2023 17:58"
$ws.Range("N109").Value = 3.54
$ws.Range("O109").Value = "30/10/2023 18:12"
$ws.Range("P109").Value = 3.33
$ws.Range("Q109").Value = "06/11/2023 17:52"
$ws.Range("R109").Value = 3.08
$ws.Range("S109").Value = "30/10/2023 18:12"
$ws.Range("T109").Value = 3.66
$ws.Range("U109").Value = "06/11/2023 17:58"
$ws.Range("V109").Value = "https://www.betexplorer.com/football/turkey/super-lig/gaziantep-rizespor/pUkr0Zc8/"

$ws.Range("A110").Value = 109
$ws.Range("B110").Value = "turkey"
$ws.Range("C110").Value = "super-lig"
$ws.Range("D110").Value = "2023-2024"
$ws.Range("E110").Value = 45236.75
$ws.Range("F110").Value = "Kayserispor"
$ws.Range("G110").Value = 1
$ws.Range("H110").Value = "Alanyaspor"
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 2.09
$ws.Range("K110").Value = "30/10/2023 02:12"
$ws.Range("L110").Value = 1.96
$ws.Range("M110").Value = "06/11/2023 17:58"
$ws.Range("N110").Value = 3.7
$ws.Range("O110").Value = "30/10/2023 02:12"
$ws.Range("P110").Value = 3.6
$ws.Range("Q110").Value = "06/11/2023 17:59"
$ws.Range("R110").Value = 3.45
$ws.Range("S110").Value = "30/10/2023 02:12"
$ws.Range("T110").Value = 4.15
$ws.Range("U110").Value = "06/11/2023 17:58"
$ws.Range("V110").Value = "https://www.betexplorer.com/football/turkey/super-lig/kayserispor-alanyaspor/YXnevhyk/"

Write-Output "Applied swaps/rotation and appended rows 109-110."
